$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.639.62'
$ws.Range("E2").Value = '  +3.27%  '
$ws.Range("D3").Value = '2.319.00'
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.43'
$ws.Range("E5").Value = '  +1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.31'
$ws.Range("E6").Value = '  +6.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.538'
$ws.Range("E7").Value = '  +2.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("E9").Value = '  +8.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.07'
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("E11").Value = '  +3.74%  '
$ws.Range("E12").Value = '  -0.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.06'
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("D14").Value = '2.682.38'
$ws.Range("E14").Value = '  +2.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.08'
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").Value = '2.321.48'
$ws.Range("E16").Value = '  +2.57%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.815'
$ws.Range("E17").Value = '  +2.52%  '
$ws.Range("D18").Value = '43.545.04'
$ws.Range("E18").Value = '  +3.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.53'
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("D20").Value = '0.0₃0931'
$ws.Range("E20").Value = '  +2.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.17'
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.48'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.59'
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("E24").Value = '  +6.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.66'
$ws.Range("E25").Value = '  +3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.00'
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.75'
$ws.Range("E28").Value = '  +4.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.54'
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.69'
$ws.Range("E31").Value = '  +1.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '167.96'
$ws.Range("E32").Value = '  +3.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.33'
$ws.Range("E33").Value = '  +1.76%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.53'
$ws.Range("E35").Value = '  +6.87%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.12'
$ws.Range("E36").Value = '  -1.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0748'
$ws.Range("E37").Value = '  +1.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.70'
$ws.Range("E38").Value = '  +2.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.107'
$ws.Range("E39").Value = '  +2.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.87'
$ws.Range("E40").Value = '  +2.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.34'
$ws.Range("E42").Value = '  +6.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.25'
$ws.Range("E43").Value = '  +7.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.31'
$ws.Range("E45").Value = '  +3.66%  '
$ws.Range("D46").Value = '1.981.69'
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.02'
$ws.Range("E47").Value = '  +4.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.91'
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.83'
$ws.Range("E49").Value = '  +4.12%  '
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.57'
$ws.Range("E51").Value = '  +7.28%  '
